$wb = $excel.ActiveWorkbook

# 1) Rename the "Include from ActCode" sheet to "Include #0"
$wsInclude = $wb.Worksheets.Item("Include from ActCode")
$wsInclude.Name = "Include #0"

# 2) Update Metadata sheet: Version + Date values, and insert a new
#    "Jurisdiction" row (row 11) before "Description", pushing
#    Description/Purpose/Copyright/Immutable down by one row.
$ws1 = $wb.Worksheets.Item("Metadata")

$ws1.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"
$ws1.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# Make room for the new row by copying row14's formatting into row15
# (keeps the same cell style used throughout the table, instead of
# Excel auto-generating a brand-new style entry).
$ws1.Range("A14:B14").Copy()
$ws1.Range("A15:B15").PasteSpecial(-4122)

# Shift rows 11-14 down into rows 12-15 (bottom-up to avoid clobbering).
$ws1.Range("A15").Value = $ws1.Range("A14").Value2
$ws1.Range("B15").Value = $ws1.Range("B14").Value2

$ws1.Range("A14").Value = $ws1.Range("A13").Value2
$ws1.Range("B14").Value = $ws1.Range("B13").Value2

$ws1.Range("A13").Value = $ws1.Range("A12").Value2
$ws1.Range("B13").Value = $ws1.Range("B12").Value2

$ws1.Range("A12").Value = $ws1.Range("A11").Value2
$ws1.Range("B12").Value = $ws1.Range("B11").Value2

# Row 11 becomes the new "Jurisdiction" property with an empty value.
$ws1.Range("A11").Value = "Jurisdiction"
$ws1.Range("B11").Value = ""
